$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.201.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.31%  "
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'2.027.67"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.12%  "
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Value = "'  +0.02%  "
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'226.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.24%  "
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'0.609"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.89%  "
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Value = "'  +0.04%  "
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'55.19"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -3.07%  "
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').Value = "'  -1.70%  "
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'0.0786"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.06%  "
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').Value = "'  -5.24%  "
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = "'2.321.67"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.42%  "
$ws.Range('E12').Style = 'Normal'

$ws.Range('E13').Value = "'  -4.23%  "
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'20.30"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -2.67%  "
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'0.744"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -1.66%  "
$ws.Range('E15').Style = 'Normal'

$ws.Range('E16').Value = "'  -2.06%  "
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'2.025.34"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.53%  "
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'37.166.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.06%  "
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'6.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +5.92%  "
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'68.96"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.78%  "
$ws.Range('E20').Style = 'Normal'

$ws.Range('E21').Value = "'  -1.22%  "
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'224.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.34%  "
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').Value = "'  +0.10%  "
$ws.Range('E23').Style = 'Normal'

$ws.Range('E24').Value = "'  +1.86%  "
$ws.Range('E24').Style = 'Normal'

$ws.Range('E25').Value = "'  -5.22%  "
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').Value = "'Monero"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'166.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.07%  "
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').Value = "'Cosmos"
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = "'9.27"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -4.89%  "
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'0.127"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.28%  "
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'18.75"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.63%  "
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').Value = "'  -2.18%  "
$ws.Range('E30').Style = 'Normal'

$ws.Range('E31').Value = "'  -1.48%  "
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').Value = "'  -0.41%  "
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'0.0614"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.75%  "
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'4.47"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -2.87%  "
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Value = "'  -4.84%  "
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = "'1.87"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.80%  "
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').Value = "'  +0.12%  "
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').Value = "'  +6.29%  "
$ws.Range('E38').Style = 'Normal'

$ws.Range('E39').Value = "'  -4.36%  "
$ws.Range('E39').Style = 'Normal'

$ws.Range('B40').Value = "'Maker"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'1.473.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -1.40%  "
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = "'VeChain"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'0.0216"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.34%  "
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').Value = "'95.96"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.85%  "
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = "'16.39"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -4.45%  "
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'0.0911"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -3.59%  "
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Value = "'  -2.25%  "
$ws.Range('E45').Style = 'Normal'

$ws.Range('E46').Value = "'  -5.28%  "
$ws.Range('E46').Style = 'Normal'

$ws.Range('E47').Value = "'  +2.28%  "
$ws.Range('E47').Style = 'Normal'

$ws.Range('E49').Value = "'  +0.35%  "
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').Value = "'FTXToken"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'3.66"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -7.29%  "
$ws.Range('E50').Style = 'Normal'

$ws.Range('B51').Value = "'RocketPoolETH"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.207.24"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -1.42%  "
$ws.Range('E51').Style = 'Normal'
